$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    under the title heading.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# 2. Insert a new paragraph right before the final "Prompt: ..." paragraph
#    containing a bold run with the page title.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$insertionRange = $lastPara.Range
$insertionRange.Collapse(1)            # wdCollapseStart
$insertionRange.InsertParagraphBefore()

$newPara = $d.Paragraphs($n)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Exotic Cats Slot Game for Free - Review 2021</w:t></w:r></w:p>'
$newPara.Range.InsertXML($titleXml) | Out-Null

# 3. Replace the text of the final "Prompt: ..." paragraph with the old
#    meta description text (the run keeps its existing italic formatting).
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.MoveEnd(1, -1) | Out-Null   # wdCharacter, exclude the paragraph mark
$finalRange.Text = "Explore exotic jungles and win big! Our review covers the features, pros, and cons of Exotic Cats online slot game. Play now for free with exciting bonuses!"
